$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update product data row (row 2) to reflect the new "Sony" product import values
$ws.Range("E2").Value = "23000"
$ws.Range("O2").Value = "Sony"
$ws.Range("F2").Value = "22000"

# Clear out the columns that are no longer populated for this product
$ws.Range("J2:K2").Clear()
$ws.Range("P2").Clear()
$ws.Range("S2:AH2").ClearContents()

# Move the active selection to match the saved view state
$ws.Range("F3").Select()
